$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("chatobogan") stat updates
$ws.Range("D7").Value = 41
$ws.Range("F7").Value = 21.49466666666666
$ws.Range("G7").Value = 2085
$ws.Range("H7").Value = 986
$ws.Range("I7").Value = 288
$ws.Range("J7").Value = 272
$ws.Range("K7").Value = 4278
$ws.Range("L7").Value = 221
$ws.Range("M7").Value = 160
$ws.Range("N7").Value = 475
$ws.Range("O7").Value = 5.390243902439025
$ws.Range("P7").Value = 3.902439024390244
$ws.Range("Q7").Value = 11.58536585365854
$ws.Range("R7").Value = 50.85
$ws.Range("S7").Value = 31.46
